$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '62.583.67'
$ws.Range("E2").Value = '  +2.77%  '

$ws.Range("D3").Value = '2.951.60'
$ws.Range("E3").Value = '  +2.56%  '

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '0.999'
$ws.Range("E4").Value = '  -0.03%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '589.09'
$ws.Range("E5").Value = '  +0.18%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '146.51'
$ws.Range("E6").Value = '  +5.34%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.999'
$ws.Range("E7").Value = '  -0.17%  '

$ws.Range("D8").Value = '2.953.90'
$ws.Range("E8").Value = '  +2.62%  '

$ws.Range("E9").Value = '  +3.31%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '6.96'
$ws.Range("E10").Value = '  +2.60%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.150'
$ws.Range("E11").Value = '  +9.69%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.436'
$ws.Range("E12").Value = '  +1.49%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.0000234'
$ws.Range("E13").Value = '  +7.45%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '32.30'
$ws.Range("E14").Value = '  +0.13%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.126'
$ws.Range("E15").Value = '  -0.89%  '

$ws.Range("D16").Value = '3.437.01'
$ws.Range("E16").Value = '  +2.62%  '

$ws.Range("D17").Value = '62.538.75'
$ws.Range("E17").Value = '  +2.68%  '

$ws.Range("B18").Value = 'WrappedEther'
$ws.Range("C18").Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range("D18").Value = '2.951.44'
$ws.Range("E18").Value = '  +2.60%  '

$ws.Range("B19").Value = 'Polkadot'
$ws.Range("C19").Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '6.65'
$ws.Range("E19").Value = '  +2.38%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '435.16'
$ws.Range("E20").Value = '  +2.16%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '13.44'
$ws.Range("E21").Value = '  +1.50%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '0.663'
$ws.Range("E22").Value = '  +1.67%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '6.97'
$ws.Range("E23").Value = '  +0.91%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '11.15'
$ws.Range("E24").Value = '  +7.21%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '80.06'
$ws.Range("E25").Value = '  +0.58%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '11.90'
$ws.Range("E26").Value = '  +4.85%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '2.10'
$ws.Range("E27").Value = '  +3.02%  '

$ws.Range("E28").Value = '  +0.04%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '7.22'
$ws.Range("E29").Value = '  +7.89%  '

$ws.Range("E30").Value = '  +4.93%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '2.58'
$ws.Range("E31").Value = '  +2.18%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '0.0000102'
$ws.Range("E32").Value = '  +20.13%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.109'
$ws.Range("E33").Value = '  +3.96%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '26.22'
$ws.Range("E34").Value = '  +2.38%  '

$ws.Range("E35").Value = '  -0.08%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.992'
$ws.Range("E36").Value = '  +3.08%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '5.58'
$ws.Range("E37").Value = '  +2.71%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '3.03'
$ws.Range("E38").Value = '  +8.78%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '49.62'
$ws.Range("E39").Value = '  +1.35%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '2.02'
$ws.Range("E40").Value = '  +6.46%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '8.36'
$ws.Range("E41").Value = '  +0.18%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.115'
$ws.Range("E42").Value = '  -1.79%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.274'
$ws.Range("E43").Value = '  +4.12%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '39.25'
$ws.Range("E44").Value = '  +1.34%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '135.38'
$ws.Range("E45").Value = '  +2.21%  '

$ws.Range("D46").Value = '2.678.84'
$ws.Range("E46").Value = '  +0.76%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.0336'
$ws.Range("E47").Value = '  +1.33%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '356.17'
$ws.Range("E48").Value = '  +3.90%  '

$ws.Range("E50").Value = '  +2.32%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '22.68'
$ws.Range("E51").Value = '  +0.95%  '
